$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = 0.9909978336781939
$ws.Range("D2").Value = 0.9859273163094625
$ws.Range("E2").Value = 0.9836697137479801
$ws.Range("F2").Value = 0.9820982118971932
$ws.Range("G2").Value = 0.9814122514125772
$ws.Range("H2").Value = 0.9831167317301314
$ws.Range("I2").Value = 0.982985637787423
$ws.Range("J2").Value = 0.9812134968852803
$ws.Range("K2").Value = 0.9807265980869618

# Row 3 updates
$ws.Range("A3").Value = 30
$ws.Range("C3").Value = 0.990662199492934
$ws.Range("D3").Value = 0.9851223533008366
$ws.Range("E3").Value = 0.9821868714552735
$ws.Range("F3").Value = 0.9799379429934298
$ws.Range("G3").Value = 0.9792504638022622
$ws.Range("H3").Value = 0.9825197307608055
$ws.Range("I3").Value = 0.9815017607672354
$ws.Range("J3").Value = 0.9783770057337925
$ws.Range("K3").Value = 0.9785632939529811
